# Commit: "Code committed by arjun on 12-10-2019 by arjun"
# Adds a new API entry (v1-flows eligibility endpoint) to the API_CONSOLE_INFO
# sheet, records the corresponding ticket numbers (ONREG-xxxxx) against the
# existing INPUT_SHEET test rows, appends a new test row, and leaves
# INPUT_SHEET as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("INPUT_SHEET")
$wsApi   = $wb.Worksheets.Item("API_CONSOLE_INFO")

# --- INPUT_SHEET -----------------------------------------------------------

# Extend the header row with a new (blank) column K, matching the existing
# blank text cells E1:J1.
$wsInput.Range("K1").Value = "'"
$wsInput.Range("K1").Style = "Normal"

# Record the ONREG ticket associated with each existing BAN/SM_USER row.
$wsInput.Range("B2").Value = "ONREG-19126"
$wsInput.Range("B3").Value = "ONREG-19127"
$wsInput.Range("B4").Value = "ONREG-17643"
$wsInput.Range("B5").Value = "ONREG-17765"
$wsInput.Range("B6").Value = "ONREG-19129"
$wsInput.Range("B7").Value = "ONREG-19067"
$wsInput.Range("B8").Value = "ONREG-19039"
$wsInput.Range("B9").Value = "ONREG-18564"

# Append the new test row.
$wsInput.Range("A10").Value = "434349902"
$wsInput.Range("B10").Value = "AugUser040"
$wsInput.Range("C10").Value = "NA"

# --- API_CONSOLE_INFO --------------------------------------------------

# Add the new "v1-flows" eligibility API entry.
$wsApi.Range("A9").Value = "v1-flows"
$wsApi.Range("B9").Value = 'https://st1-apiservices-sen.test.sprint.com:8441/api/process/eligibility/v1/flows?type=CHANGE_SERVICES&subscriberId=$SUBSCRIBER&accountId=$BAN&accountSubscriberIds=$SUBSCRIBER&role=ACCOUNT_OWNER'
$wsApi.Range("C9").Value = "GET"

# --- Selection / active sheet state ----------------------------------------

$wsApi.Range("A10").Select()
$wsInput.Activate()
$wsInput.Range("C10").Select()
